# Generate Report for Handoff
# A new handoff was generated for the e2e test markdown file (new GUID
# 841dada6-e132-4bf1-acfa-efa3f74ee009 replacing a5cfe79f-3c5f-4df4-8006-b4782083b8ec),
# with a fresh hash (13b7a6bd80b9b66a5ebac852346bef08ff6134ba) and handoff
# timestamps, and the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns reset because the new handoff has not
# been handed back yet.

$wb = $excel.ActiveWorkbook

$oldGuid = "a5cfe79f-3c5f-4df4-8006-b4782083b8ec"
$newGuid = "841dada6-e132-4bf1-acfa-efa3f74ee009"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Range("B2").Value = "e2e\$newGuid.md"
$ws1.Range("G2").Value = "2016-09-07 09:32:59"

# Refresh the hyperlink display text on B2, keeping the same target URL.
$b2 = $ws1.Range("B2")
$oldTargetOverview = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a02f2c01a76ca446607a0ef475c3500a7439b9bc/e2e/$oldGuid.md"
$b2.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($b2, $oldTargetOverview, "", "", "e2e\$newGuid.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "$newGuid.md"
$ws2.Range("G2").Value = "$newGuid.13b7a6bd80b9b66a5ebac852346bef08ff6134ba.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-07 09:32:53"
$ws2.Range("I2").Hyperlinks.Delete()
$ws2.Range("I2").Value = ""
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

$a2zh = $ws2.Range("A2")
$oldTargetZh = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a02f2c01a76ca446607a0ef475c3500a7439b9bc/e2e/$oldGuid.md"
$a2zh.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($a2zh, $oldTargetZh, "", "", "$newGuid.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "$newGuid.md"
$ws3.Range("G2").Value = "$newGuid.13b7a6bd80b9b66a5ebac852346bef08ff6134ba.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-07 09:32:59"
$ws3.Range("I2").Hyperlinks.Delete()
$ws3.Range("I2").Value = ""
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$a2de = $ws3.Range("A2")
$oldTargetDe = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a02f2c01a76ca446607a0ef475c3500a7439b9bc/e2e/$oldGuid.md"
$a2de.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($a2de, $oldTargetDe, "", "", "$newGuid.md")

# ---------------------------------------------------------------------------
# Column width tweaks on the "Latest Target File" / "Latest Handback File"
# columns (I/J) for the zh-cn and de-de sheets, now that they hold short
# empty values instead of long file names.
# ---------------------------------------------------------------------------
$ws2.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws2.Columns.Item(10).ColumnWidth = 21.7054770333426
$ws3.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws3.Columns.Item(10).ColumnWidth = 21.7054770333426
